{"js": "// Insert \" repeatedly\" immediately before the final period of the sentence:\n// \"...export the dataframes onto the database.\" ->\n// \"...export the dataframes onto the database repeatedly.\"\n//\n// Strategy: locate the unique sentence fragment, split it down to the word\n// \"database\" (without its trailing period) using Range.split(), then insert\n// \" repeatedly\" at the end of that sub-range (i.e. right before the period).\n\nconst body = context.document.body;\n\nconst anchorText = \"the dataframes onto the database.\";\nconst searchResults = body.search(anchorText, { matchCase: false, matchWholeWord: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not locate target sentence fragment.\");\n}\n\nconst matchRange = searchResults.items[0];\n\n// Split into words (keep trailing space attached to each word, last item keeps \".\")\nconst wordRanges = matchRange.split([\" \"], false, true, false);\nwordRanges.load(\"text\");\nawait context.sync();\n\nconst lastWordRange = wordRanges.items[wordRanges.items.length - 1]; // \"database.\"\n\n// Split off the trailing \".\" so we have a sub-range ending right after \"database\".\nconst noPeriodRanges = lastWordRange.split([\".\"], false, true, false);\nnoPeriodRanges.load(\"text\");\nawait context.sync();\n\nconst databaseRange = noPeriodRanges.items[0]; // \"database\" (no period)\n\n// Collapsed insertion point right after \"database\" and before the period.\nconst insertionPoint = databaseRange.getRange(\"End\");\ninsertionPoint.insertText(\" repeatedly\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Insert \" repeatedly\" immediately before the final period of the sentence:\n# \"...export the dataframes onto the database.\" ->\n# \"...export the dataframes onto the database repeatedly.\"\n\n$d = $word.ActiveDocument\n\n# Locate the unique trailing fragment of the sentence (includes the final period).\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \"onto the database.\"\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not locate target sentence fragment.\"\n}\n\n# $searchRange now spans the matched text (\"onto the database.\").\n# Build a collapsed insertion-point range right before the trailing period\n# (i.e. one character before the end of the match) and insert the addition there.\n$insertionRange = $d.Range($searchRange.End - 1, $searchRange.End - 1)\n$insertionRange.InsertBefore(\" repeatedly\")\n"}
